{"js": "// A\u00f1adidos los colores base para la p\u00e1gina\n// Insert 5 new paragraphs right after the \"Letra \u2192 #05161A\" paragraph:\n//   (empty), \"Alternativos\", \"#0C7075\", \"#0F969C\", \"#294D61\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph whose text is exactly \"Letra \u2192 #05161A\".\nlet anchor = null;\nfor (const p of paragraphs.items) {\n  if (p.text === \"Letra \u2192 #05161A\") {\n    anchor = p;\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error('Could not find paragraph \"Letra \u2192 #05161A\"');\n}\n\nconst newTexts = [\"\", \"Alternativos\", \"#0C7075\", \"#0F969C\", \"#294D61\"];\n\nlet current = anchor;\nfor (const text of newTexts) {\n  current = current.insertParagraph(text, Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "# A\u00f1adidos los colores base para la p\u00e1gina\n# Insert 5 new paragraphs right after the \"Letra \u2192 #05161A\" paragraph:\n#   (empty), \"Alternativos\", \"#0C7075\", \"#0F969C\", \"#294D61\"\n\n$d = $word.ActiveDocument\n\n$searchRange = $d.Content\n$found = $searchRange.Find.Execute(\"Letra \u2192 #05161A\")\nif (-not $found) {\n    throw 'Could not find paragraph \"Letra \u2192 #05161A\"'\n}\n\n$anchor = $searchRange.Paragraphs(1)\n\n$lines = @(\"\", \"Alternativos\", \"#0C7075\", \"#0F969C\", \"#294D61\")\n\nforeach ($line in $lines) {\n    $insertionPoint = $anchor.Range.Duplicate\n    $insertionPoint.Collapse(0)   # wdCollapseEnd\n    $insertionPoint.InsertParagraphAfter()\n\n    $newPara = $anchor.Next()\n    if ($line -ne \"\") {\n        $newPara.Range.Text = $line\n    }\n\n    $anchor = $newPara\n}\n\n$d.Save()\n"}
